$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each updated cell as literal text (NumberFormat "@" forces text entry,
# then resetting the Style back to "Normal" avoids leaving a stray number format
# on the cell so the final style matches the original, unstyled cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '62.373.26'
Set-TextValue $ws.Range('E2') '  -2.32%  '
Set-TextValue $ws.Range('D3') '2.640.96'
Set-TextValue $ws.Range('E3') '  -3.47%  '
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '549.33'
Set-TextValue $ws.Range('E5') '  -2.93%  '
Set-TextValue $ws.Range('D6') '154.09'
Set-TextValue $ws.Range('E6') '  -4.64%  '
Set-TextValue $ws.Range('E7') '  +0.09%  '
Set-TextValue $ws.Range('E8') '  -1.76%  '
Set-TextValue $ws.Range('E9') '  -4.84%  '
Set-TextValue $ws.Range('E10') '  -4.47%  '
Set-TextValue $ws.Range('D11') '5.43'
Set-TextValue $ws.Range('E11') '  -3.69%  '
Set-TextValue $ws.Range('D12') '0.361'
Set-TextValue $ws.Range('E12') '  -5.07%  '
Set-TextValue $ws.Range('D13') '3.106.80'
Set-TextValue $ws.Range('E13') '  -3.58%  '
Set-TextValue $ws.Range('D14') '25.64'
Set-TextValue $ws.Range('E14') '  -5.08%  '
Set-TextValue $ws.Range('D15') '62.290.28'
Set-TextValue $ws.Range('D16') '0.0000143'
Set-TextValue $ws.Range('E16') '  -4.41%  '
Set-TextValue $ws.Range('D17') '2.644.85'
Set-TextValue $ws.Range('E17') '  -3.54%  '
Set-TextValue $ws.Range('D19') '4.52'
Set-TextValue $ws.Range('E19') '  -4.56%  '
Set-TextValue $ws.Range('D20') '339.55'
Set-TextValue $ws.Range('E20') '  -4.48%  '
Set-TextValue $ws.Range('E21') '  -8.04%  '
Set-TextValue $ws.Range('D22') '0.997'
Set-TextValue $ws.Range('E22') '  -0.20%  '
Set-TextValue $ws.Range('D23') '0.501'
Set-TextValue $ws.Range('E23') '  -3.90%  '
Set-TextValue $ws.Range('E24') '  -3.13%  '
Set-TextValue $ws.Range('E25') '  -1.39%  '
Set-TextValue $ws.Range('E26') '  +0.02%  '
Set-TextValue $ws.Range('D27') '7.99'
Set-TextValue $ws.Range('E27') '  -4.87%  '
Set-TextValue $ws.Range('B28') 'Fetch.AI'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D28') '1.36'
Set-TextValue $ws.Range('E28') '  +0.90%  '
Set-TextValue $ws.Range('B29') 'PEPE'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D29') '0.0₃0833'
Set-TextValue $ws.Range('E29') '  -8.40%  '
Set-TextValue $ws.Range('D30') '7.06'
Set-TextValue $ws.Range('E30') '  -1.48%  '
Set-TextValue $ws.Range('E31') '  -5.69%  '
Set-TextValue $ws.Range('D32') '160.21'
Set-TextValue $ws.Range('E32') '  -3.73%  '
Set-TextValue $ws.Range('E33') '  +0.00%  '
Set-TextValue $ws.Range('E34') '  -3.84%  '
Set-TextValue $ws.Range('E35') '  -3.61%  '
Set-TextValue $ws.Range('D36') '19.18'
Set-TextValue $ws.Range('E36') '  -4.74%  '
Set-TextValue $ws.Range('E37') '  -4.91%  '
Set-TextValue $ws.Range('D38') '334.87'
Set-TextValue $ws.Range('E38') '  -2.92%  '
Set-TextValue $ws.Range('D39') '6.09'
Set-TextValue $ws.Range('E39') '  -3.35%  '
Set-TextValue $ws.Range('D40') '0.903'
Set-TextValue $ws.Range('E40') '  -7.68%  '
Set-TextValue $ws.Range('D41') '38.01'
Set-TextValue $ws.Range('E41') '  -2.19%  '
Set-TextValue $ws.Range('D42') '3.92'
Set-TextValue $ws.Range('E42') '  -4.30%  '
Set-TextValue $ws.Range('D43') '0.999'
Set-TextValue $ws.Range('E43') '  +0.10%  '
Set-TextValue $ws.Range('D44') '20.38'
Set-TextValue $ws.Range('E44') '  -6.36%  '
Set-TextValue $ws.Range('D45') '0.608'
Set-TextValue $ws.Range('E45') '  -3.67%  '
Set-TextValue $ws.Range('D46') '19.77'
Set-TextValue $ws.Range('E46') '  -6.63%  '
Set-TextValue $ws.Range('E47') '  -0.75%  '
Set-TextValue $ws.Range('D48') '0.0546'
Set-TextValue $ws.Range('E48') '  -7.06%  '
Set-TextValue $ws.Range('D49') '0.0959'
Set-TextValue $ws.Range('E49') '  -4.09%  '
Set-TextValue $ws.Range('D50') '127.17'
Set-TextValue $ws.Range('E50') '  -3.99%  '
Set-TextValue $ws.Range('D51') '0.0237'
Set-TextValue $ws.Range('E51') '  -5.77%  '
